# Applies the commit "remove If nodes, fix bugs, remove word instructions"
# to the Circularity Tool Excel report template.
#
# Summary of the functional changes:
#  - MCI_results: the two literal "origin material" fractions that used to
#    sum with the other five columns to 1.0 are replaced by raw weight-like
#    numbers (the previous 0.4/0.6 split no longer works as a percentage).
#  - Auxiliary/Report: the two averages that read those columns are fixed so
#    they again produce a 0-1 ratio by dividing by the average element
#    weight, and the two recomputed Auxiliary cells are bolded to flag them.
#  - The Auxiliary helper sheet is hidden from the tab strip.
#  - A few leftover cell selections/scroll positions are nudged (cosmetic).

$wb = $excel.ActiveWorkbook

# --- MCI_results: bugfixed raw inputs (was a 0-1 fraction, now a weight-like
#     quantity) for the two elements in the demo table -----------------------
$mci = $wb.Worksheets.Item("MCI_results")
$mci.Range("I2").Value = 800
$mci.Range("O2").Value = 500
$mci.Range("I3").Value = 100
$mci.Range("O3").Value = 200

$mci.Activate()
$mci.Range("O4").Select()

# --- Auxiliary: normalise the two averages by the average weight, and bold
#     the two cells to highlight the fixed formulas --------------------------
$aux = $wb.Worksheets.Item("Auxiliary")
$aux.Range("C2").Formula = "=AVERAGE(MCI[Virgin material (V)])/AVERAGE(MCI[Weight])"
$aux.Range("C2").Font.Bold = $true
$aux.Range("F4").Formula = "=AVERAGE(MCI[Linear waste (W0)])/AVERAGE(MCI[Weight])"
$aux.Range("F4").Font.Bold = $true

$aux.Activate()
$aux.Range("F3").Select()

# hide the Auxiliary helper sheet now that it is no longer meant to be seen
$aux.Visible = $false

# --- Report: same formula fix as Auxiliary, same cached result -------------
$report = $wb.Worksheets.Item("Report")
$report.Range("E21").Formula = "=AVERAGE(MCI[Virgin material (V)])/AVERAGE(MCI[Weight])"
$report.Range("E27").Formula = "=AVERAGE(MCI[Linear waste (W0)])/AVERAGE(MCI[Weight])"

$report.Activate()
$report.Range("D23").Select()
